# Updated symbol list on Fri Dec 23 07:22:15 UTC 2022 with GitHub Actions
#
# All cells in this sheet are stored as plain text (General format, no
# number formatting) even though several columns hold numeric-looking
# strings (prices, etc). A naive `$cell.Value = "1.23"` assignment lets
# Excel auto-coerce the text into a real Number, which would change the
# cell's stored type/format and not match the source data. To keep the
# values as genuine text we prefix numeric-looking strings with a literal
# apostrophe (Excel's "treat as text" marker) and then restore the cell's
# original ("Normal") style afterwards, since the apostrophe entry also
# marks the cell with a quote-prefix style that the source file doesn't
# have.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Sheet, [string]$A1, [string]$Text)

    $cell = $Sheet.Range($A1)
    if ($Text -match '^-?\d+(\.\d+)?$') {
        # Numeric-looking value: force text entry via leading apostrophe,
        # then reset the style so we don't leave a stray quote-prefix
        # format behind (keeps styling identical to the source cell).
        $cell.Value = "'" + $Text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $Text
    }
}

# --- Price (column D) updates ---
Set-CellText $ws "D2"  "246.73"
Set-CellText $ws "D3"  "21.92"
Set-CellText $ws "D4"  "5.416"
Set-CellText $ws "D5"  "0.05785"
Set-CellText $ws "D7"  "6.336"
Set-CellText $ws "D8"  "0.8075"
Set-CellText $ws "D9"  "0.9459"
Set-CellText $ws "D10" "0.1427"
Set-CellText $ws "D11" "0.07481"
Set-CellText $ws "D14" "4.153"
Set-CellText $ws "D15" "0.09412"
Set-CellText $ws "D16" "0.001589"
Set-CellText $ws "D17" "0.04820"
Set-CellText $ws "D18" "0.0005889"
Set-CellText $ws "D21" "0.0009928"
Set-CellText $ws "D23" "3.773"
Set-CellText $ws "D24" "2.233"
Set-CellText $ws "D25" "0.3233"
Set-CellText $ws "D40" "0.03895"

# --- Rows 41-43: coin list rotated (KickToken/BKEXToken/CEJI reshuffled) ---
Set-CellText $ws "B41" "BKEXToken"
Set-CellText $ws "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-CellText $ws "D41" "0.1075"
Set-CellText $ws "E41" "40BKEXTokenBKK"

Set-CellText $ws "B42" "CEJI"
Set-CellText $ws "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-CellText $ws "D42" "0.002660"
Set-CellText $ws "E42" "41CEJICEJI"

Set-CellText $ws "B43" "KickToken"
Set-CellText $ws "C43" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-CellText $ws "D43" "0.006359"
Set-CellText $ws "E43" "42KickTokenKICK"

# --- remaining price (column D) updates ---
Set-CellText $ws "D44" "0.006331"
Set-CellText $ws "D45" "0.00005591"
Set-CellText $ws "D48" "0.1438"
Set-CellText $ws "D49" "0.00002100"

Write-Output "edit.ps1 applied 35 cell updates"
